$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns before column F (shift existing F:U right to K:Z),
# so the new "meta" block lands in F:J.
$ws.Columns("F:J").Insert()

# New header labels (row 1) for the inserted columns.
$ws.Range("F1").Value = "meta"
$ws.Range("G1").Value = "meta_avg"
$ws.Range("H1").Value = "meta_std"
$ws.Range("I1").Value = "meta_min"
$ws.Range("J1").Value = "meta_max"

# Apply the same number format ("R$ #,##0.00") used by the neighboring
# "arrecadado_*" columns to the new data cells (rows 2-4 only, so the
# header row keeps its original bold/border/centered style).
$ws.Range("F2:J4").NumberFormat = "R$ #,##0.00"

# Row 2 (aon) new "meta" values
$ws.Range("F2").Value = 13973042.60019265
$ws.Range("G2").Value = 16834.99108456945
$ws.Range("H2").Value = 17015.69760983049
$ws.Range("I2").Value = 31.89582864100442
$ws.Range("J2").Value = 189313.7035611726

# Row 3 (flex) new "meta" values
$ws.Range("F3").Value = 15599716.7029188
$ws.Range("G3").Value = 11279.62162177787
$ws.Range("H3").Value = 16430.30708090436
$ws.Range("I3").Value = 12.04441558726698
$ws.Range("J3").Value = 198811.9434626772

# Row 4 (sub) new "meta" values
$ws.Range("F4").Value = 165199.0578149446
$ws.Range("G4").Value = 1205.832538795216
$ws.Range("H4").Value = 2163.288658625352
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 21176.91783511972

# Negligible last-digit re-serialization of pre-existing (untouched) values
# carried over from the shift; restate them exactly to match the source.
$ws.Range("M4").Value = 650.5808076401025
$ws.Range("V2").Value = 423.0192251466749
$ws.Range("P4").Value = 21.28348419490776
